# [docx writer] Make list para properties go first
#
# The numPr element in each list paragraph's pPr currently serializes
# numId before ilvl. Re-assert each list paragraph's list level so the
# writer re-emits numPr with ilvl first, numId second (matches what
# Word Online expects).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.ListFormat.ListType -ne 0) {
        $p.Range.ListFormat.ListLevelNumber = $p.Range.ListFormat.ListLevelNumber
    }
}
